# Generate Report for Handback
#
# The 8c1d74b3-... file has now been handed back (in sync with en-US).
# It moves to the top of each table (row 2), gets a new status, and gets
# its "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated. The 6e73dd2d-... file (still awaiting
# handoff) drops to row 3, unchanged otherwise.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# URLs reused / needed for hyperlinks
# ---------------------------------------------------------------------
$urlMd8c1d   = "https://github.com/OpenLocalizationTest/oltest/blob/596b1e43861d72ad231980a76b5a1d98d1ced5b4/e2e/8c1d74b3-d862-47a2-aad4-e9501adc2638.md"
$urlMd6e73   = "https://github.com/OpenLocalizationTest/oltest/blob/96a37853af021d71e9a6a338eab167b9b88022ad/e2e/6e73dd2d-a7e3-4c5a-af29-978236a9d6bd.md"
$urlLocalCfg = "https://github.com/OpenLocalizationTest/oltest/blob/596b1e43861d72ad231980a76b5a1d98d1ced5b4/.localization-config"

$urlXlf8c1dZhCn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5246fe4481a40f77775f48e9c5ad6e22246ff91b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8c1d74b3-d862-47a2-aad4-e9501adc2638.f238ee7f8e9639ffebef0130ddaec26b74c295ef.zh-cn.xlf"
$urlXlf6e73ZhCn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/502fb4aa705b27df7f36120353f9b8a617039e70/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6e73dd2d-a7e3-4c5a-af29-978236a9d6bd.5bbf06d8975a65ba18cacba0b87de39396998ad2.zh-cn.xlf"

$urlXlf8c1dDeDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/227f06d3ae556147a08c4294fe0bd2b282581fac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8c1d74b3-d862-47a2-aad4-e9501adc2638.f238ee7f8e9639ffebef0130ddaec26b74c295ef.de-de.xlf"
$urlXlf6e73DeDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ffb57ec49cb63a562bf1745cc1b83db530568e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6e73dd2d-a7e3-4c5a-af29-978236a9d6bd.5bbf06d8975a65ba18cacba0b87de39396998ad2.de-de.xlf"

$dispMd8c1d = "8c1d74b3-d862-47a2-aad4-e9501adc2638.md"
$dispMd6e73 = "6e73dd2d-a7e3-4c5a-af29-978236a9d6bd.md"
$dispLocalCfg = ".localization-config"
$dispXlf8c1dZhCn = "8c1d74b3-d862-47a2-aad4-e9501adc2638.f238ee7f8e9639ffebef0130ddaec26b74c295ef.zh-cn.xlf"
$dispXlf6e73ZhCn = "6e73dd2d-a7e3-4c5a-af29-978236a9d6bd.5bbf06d8975a65ba18cacba0b87de39396998ad2.zh-cn.xlf"
$dispXlf8c1dDeDe = "8c1d74b3-d862-47a2-aad4-e9501adc2638.f238ee7f8e9639ffebef0130ddaec26b74c295ef.de-de.xlf"
$dispXlf6e73DeDe = "6e73dd2d-a7e3-4c5a-af29-978236a9d6bd.5bbf06d8975a65ba18cacba0b87de39396998ad2.de-de.xlf"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady = "Ready for handoff"
$notLocalized = "Not to be localized"

# ---------------------------------------------------------------------
# Sheet "Overview": swap the two rows so 8c1d74b3 is listed first
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = $dispMd8c1d
$wsOverview.Range("B2").Value2 = $statusHandedBack
$wsOverview.Range("C2").Value2 = $statusHandedBack

$wsOverview.Range("A3").Value2 = $dispMd6e73
$wsOverview.Range("B3").Value2 = $statusReady
$wsOverview.Range("C3").Value2 = $statusReady

$wsOverview.Range("A4").Value2 = $dispLocalCfg
$wsOverview.Range("B4").Value2 = $notLocalized
$wsOverview.Range("C4").Value2 = $notLocalized

# Rebuild the hyperlinks on the Overview sheet in the new order
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $urlMd8c1d, [System.Type]::Missing, [System.Type]::Missing, $dispMd8c1d) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urlMd6e73, [System.Type]::Missing, [System.Type]::Missing, $dispMd6e73) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $urlLocalCfg, [System.Type]::Missing, [System.Type]::Missing, $dispLocalCfg) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 2: now the 8c1d74b3 entry - handed back, in sync
$wsZhCn.Range("A2").Value2 = $dispMd8c1d
$wsZhCn.Range("B2").Value2 = $statusHandedBack
$wsZhCn.Range("C2").Value2 = $dispXlf8c1dZhCn
$wsZhCn.Range("D2").Value2 = "2016-03-08 16:33:46"
$wsZhCn.Range("E2").Value2 = $dispMd8c1d
$wsZhCn.Range("F2").Value2 = $dispXlf8c1dZhCn
$wsZhCn.Range("G2").Value2 = "2016-03-08 16:34:15"
$wsZhCn.Range("H2").Value2 = "Include"

# Row 3: now the 6e73dd2d entry - still waiting, unchanged values
$wsZhCn.Range("A3").Value2 = $dispMd6e73
$wsZhCn.Range("B3").Value2 = $statusReady
$wsZhCn.Range("C3").Value2 = $dispXlf6e73ZhCn
$wsZhCn.Range("D3").Value2 = "2016-03-08 16:33:07"
$wsZhCn.Range("E3").ClearContents()
$wsZhCn.Range("F3").ClearContents()
$wsZhCn.Range("G3").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("H3").Value2 = "Include"

# Row 4: .localization-config, unchanged
$wsZhCn.Range("A4").Value2 = $dispLocalCfg
$wsZhCn.Range("B4").Value2 = $notLocalized
$wsZhCn.Range("D4").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("G4").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("H4").Value2 = "Ignored"

# Rebuild the hyperlinks on the zh-cn sheet in the new order
$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlMd8c1d, [System.Type]::Missing, [System.Type]::Missing, $dispMd8c1d) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $urlXlf8c1dZhCn, [System.Type]::Missing, [System.Type]::Missing, $dispXlf8c1dZhCn) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), $urlMd8c1d, [System.Type]::Missing, [System.Type]::Missing, $dispMd8c1d) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $urlXlf8c1dZhCn, [System.Type]::Missing, [System.Type]::Missing, $dispXlf8c1dZhCn) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlMd6e73, [System.Type]::Missing, [System.Type]::Missing, $dispMd6e73) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), $urlXlf6e73ZhCn, [System.Type]::Missing, [System.Type]::Missing, $dispXlf6e73ZhCn) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $urlLocalCfg, [System.Type]::Missing, [System.Type]::Missing, $dispLocalCfg) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 2: now the 8c1d74b3 entry - handed back, in sync
$wsDeDe.Range("A2").Value2 = $dispMd8c1d
$wsDeDe.Range("B2").Value2 = $statusHandedBack
$wsDeDe.Range("C2").Value2 = $dispXlf8c1dDeDe
$wsDeDe.Range("D2").Value2 = "2016-03-08 16:33:52"
$wsDeDe.Range("E2").Value2 = $dispMd8c1d
$wsDeDe.Range("F2").Value2 = $dispXlf8c1dDeDe
$wsDeDe.Range("G2").Value2 = "2016-03-08 16:34:31"
$wsDeDe.Range("H2").Value2 = "Include"

# Row 3: now the 6e73dd2d entry - still waiting, unchanged values
$wsDeDe.Range("A3").Value2 = $dispMd6e73
$wsDeDe.Range("B3").Value2 = $statusReady
$wsDeDe.Range("C3").Value2 = $dispXlf6e73DeDe
$wsDeDe.Range("D3").Value2 = "2016-03-08 16:33:13"
$wsDeDe.Range("E3").ClearContents()
$wsDeDe.Range("F3").ClearContents()
$wsDeDe.Range("G3").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("H3").Value2 = "Include"

# Row 4: .localization-config, unchanged
$wsDeDe.Range("A4").Value2 = $dispLocalCfg
$wsDeDe.Range("B4").Value2 = $notLocalized
$wsDeDe.Range("D4").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("G4").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("H4").Value2 = "Ignored"

# Rebuild the hyperlinks on the de-de sheet in the new order
$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlMd8c1d, [System.Type]::Missing, [System.Type]::Missing, $dispMd8c1d) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $urlXlf8c1dDeDe, [System.Type]::Missing, [System.Type]::Missing, $dispXlf8c1dDeDe) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), $urlMd8c1d, [System.Type]::Missing, [System.Type]::Missing, $dispMd8c1d) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $urlXlf8c1dDeDe, [System.Type]::Missing, [System.Type]::Missing, $dispXlf8c1dDeDe) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlMd6e73, [System.Type]::Missing, [System.Type]::Missing, $dispMd6e73) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), $urlXlf6e73DeDe, [System.Type]::Missing, [System.Type]::Missing, $dispXlf6e73DeDe) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $urlLocalCfg, [System.Type]::Missing, [System.Type]::Missing, $dispLocalCfg) | Out-Null
